# DeveloperGuide diagram update: LogicComponentClassDiagram.pptx
#
# 1) The auto "date last saved" footer fields (Date Placeholder shapes) on
#    the slide master and every slide layout get refreshed from
#    "7/20/17" to "4/16/2018".
# 2) The now-obsolete "UndoRedoStack" mini-diagram (a rectangle, the arrow
#    pointing into it and its "1" callout label) is removed from the
#    slide, since the app's undo/redo mechanism moved from
#    UndoRedoStack/UndoableCommand to VersionedAddressBook.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "7/20/17") {
                $shp.TextFrame.TextRange.Text = "4/16/2018"
            }
        }
    }
}

# --- 1. Refresh the date field text everywhere it appears -----------------

$master = $p.SlideMaster

# Slide master's own Date Placeholder
Update-DatePlaceholder $master.Shapes

# Every slide layout's Date Placeholder
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2. Remove the obsolete UndoRedoStack shapes from the slide -----------

$s = $p.Slides.Item(1)
$idsToDelete = @(59, 61, 63)
foreach ($targetId in $idsToDelete) {
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Id -eq $targetId) {
            $shp.Delete()
            break
        }
    }
}
